# Public school teacher age data
# Fill in the previously-blank cells in row 9 ("General elementary education")
# with the "double dagger" reporting-standards-not-met marker, matching the
# marker already used elsewhere in the sheet (e.g. K9, column K of rows 11-15).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B9").Value = "‡"
$ws.Range("D9").Value = "‡"
$ws.Range("F9").Value = "‡"
$ws.Range("J9").Value = "‡"
